# Fruta / hortaliza, semanal
# Insert a new weekly record as row 5, pushing the existing rows 5-10 down to 6-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 5 (shifts old rows 5-10 to 6-11)
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new weekly entry
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "Macroferia Regional de Talca"
$ws.Range("C5").Value = "Maule"
$ws.Range("D5").Value = "12/20/2023"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100104
$ws.Range("H5").Value = "Frutos de pepita"
$ws.Range("I5").Value = 100104004
$ws.Range("J5").Value = "Níspero"
$ws.Range("K5").Value = "Golden Nugget"
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "$/caja 10 kilos"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1500
$ws.Range("T5").Value = 10
